# Add a new "membership" worksheet (as the 3rd/last sheet) that mirrors the
# layout of Sheet1/Sheet2 (username/password + membership/amount/currency
# columns), per the "Membership details test template committed" change.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Duplicate Sheet2 to the end of the workbook so the new sheet inherits the
# same sheetPr/outlinePr, sheetFormatPr and drawing relationship conventions
# used by the existing sheets, then rename it and wipe its contents.
$ws2.Copy($null, $ws2)
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "membership"
$ws.Cells.Clear()

# Pull in the base font/formatting (style used for the header + data rows
# on Sheet1) for the whole working area before writing values.
$ws1.Range("A1:E2").Copy()
$ws.Range("A1:E2").PasteSpecial(-4122)

# Header row.
$ws.Range("A1").Value = '${username}'
$ws.Range("B1").Value = '${password}'
$ws.Range("C1").Value = '${membership}'
$ws.Range("D1").Value = '${amount}'
$ws.Range("E1").Value = '${currency}'

# Data row.
$ws.Range("A2").Value = 'Admin'
$ws.Range("B2").Value = 'admin123'
$ws.Range("C2").Value = 'ACCA'

# Amount is stored as text (Text number format), matching the source data.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '500'

$ws.Range("E2").Value = 'Indian Rupee'

# Trailing blank filler cells under "amount" (D3:D23), formatted the same
# way as the equivalent text-formatted placeholder column on Sheet2.
$ws2.Range("E1").Copy()
$ws.Range("D3:D23").PasteSpecial(-4122)
